$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.942807952859368
$ws.Range("D2").Value = 4.563323565394718
$ws.Range("E2").Value = 12.39846029903014
$ws.Range("F2").Value = 23.13584562483873
$ws.Range("G2").Value = 3.620216811455744
$ws.Range("I2").Value = 19.86122920051779
$ws.Range("K2").Value = 10.39283784207068
$ws.Range("M2").Value = 13.70291929367718
$ws.Range("O2").Value = 20.64567129158716
$ws.Range("B3").Value = 5.813417697830533
$ws.Range("D3").Value = 4.516685529333609
$ws.Range("E3").Value = 12.27229026009132
$ws.Range("F3").Value = 23.13429470482652
$ws.Range("G3").Value = 3.622226647265002
$ws.Range("I3").Value = 19.98322573029162
$ws.Range("K3").Value = 9.946387375841262
$ws.Range("M3").Value = 13.4472273269489
$ws.Range("O3").Value = 20.70308937203066
$ws.Range("B4").Value = 5.733210816678895
$ws.Range("D4").Value = 4.487456247845639
$ws.Range("E4").Value = 12.19904410531813
$ws.Range("F4").Value = 23.14090834035354
$ws.Range("G4").Value = 3.623525896958909
$ws.Range("I4").Value = 20.06190135356284
$ws.Range("K4").Value = 9.659854535169691
$ws.Range("M4").Value = 13.29075024692689
$ws.Range("O4").Value = 20.74405785119519
$ws.Range("B5").Value = 5.700379323544152
$ws.Range("D5").Value = 4.47540189567056
$ws.Range("E5").Value = 12.17029133549018
$ws.Range("F5").Value = 23.14550351547362
$ws.Range("G5").Value = 3.624071800283807
$ws.Range("I5").Value = 20.09491279055556
$ws.Range("K5").Value = 9.54006820147711
$ws.Range("M5").Value = 13.2272026445543
$ws.Range("O5").Value = 20.76218407747638
$ws.Range("B6").Value = 5.694920168356521
$ws.Range("D6").Value = 4.473391809440074
$ws.Range("E6").Value = 12.16558406603736
$ws.Range("F6").Value = 23.14638119767184
$ws.Range("G6").Value = 3.624163442120683
$ws.Range("I6").Value = 20.10045179989163
$ws.Range("K6").Value = 9.519998440805193
$ws.Range("M6").Value = 13.21666636429565
$ws.Range("O6").Value = 20.76528019969978
$ws.Range("B7").Value = 5.732768572662785
$ws.Range("D7").Value = 4.487294250372981
$ws.Range("E7").Value = 12.19865185666488
$ws.Range("F7").Value = 23.14096262364032
$ws.Range("G7").Value = 3.623533192532557
$ws.Range("I7").Value = 20.06234270534231
$ws.Range("K7").Value = 9.658251148709979
$ws.Range("M7").Value = 13.28989222227055
$ws.Range("O7").Value = 20.74429652033137
$ws.Range("B8").Value = 5.898378411588793
$ws.Range("D8").Value = 4.547369494464101
$ws.Range("E8").Value = 12.35410154227642
$ws.Range("F8").Value = 23.13373956671119
$ws.Range("G8").Value = 3.62089630195883
$ws.Range("I8").Value = 19.90251290658284
$ws.Range("K8").Value = 10.24153230823836
$ws.Range("M8").Value = 13.61470122915489
$ws.Range("O8").Value = 20.66428001024688
$ws.Range("B9").Value = 6.215227141364752
$ws.Range("D9").Value = 4.660211305420528
$ws.Range("E9").Value = 12.69082953961534
$ws.Range("F9").Value = 23.17964520948732
$ws.Range("G9").Value = 3.616240305454796
$ws.Range("I9").Value = 19.61887147007902
$ws.Range("K9").Value = 11.28324198282091
$ws.Range("M9").Value = 14.25198978089723
$ws.Range("O9").Value = 20.55293211985332
$ws.Range("B10").Value = 6.44082574236813
$ws.Range("D10").Value = 4.739754198062028
$ws.Range("E10").Value = 12.95534782812
$ws.Range("F10").Value = 23.24994013175752
$ws.Range("G10").Value = 3.61313008363188
$ws.Range("I10").Value = 19.4284634253308
$ws.Range("K10").Value = 11.98234770527418
$ws.Range("M10").Value = 14.71550588476573
$ws.Range("O10").Value = 20.49919311091648
$ws.Range("B11").Value = 6.541439636274565
$ws.Range("D11").Value = 4.775142035442228
$ws.Range("E11").Value = 13.07887501852657
$ws.Range("F11").Value = 23.28981306128655
$ws.Range("G11").Value = 3.611781874055985
$ws.Range("I11").Value = 19.34570975961178
$ws.Range("K11").Value = 12.28537393349288
$ws.Range("M11").Value = 14.92436975555924
$ws.Range("O11").Value = 20.48089423224153
$ws.Range("B12").Value = 6.57921548753408
$ws.Range("D12").Value = 4.788422366557443
$ws.Range("E12").Value = 13.12606555628429
$ws.Range("F12").Value = 23.30604051184133
$ws.Range("G12").Value = 3.611280871188959
$ws.Range("I12").Value = 19.31492597703357
$ws.Range("K12").Value = 12.3979235679288
$ws.Range("M12").Value = 15.00309770016975
$ws.Range("O12").Value = 20.47485252565985
$ws.Range("B13").Value = 6.57109474437704
$ws.Range("D13").Value = 4.785567648121837
$ws.Range("E13").Value = 13.11588458074831
$ws.Range("F13").Value = 23.30249558084729
$ws.Range("G13").Value = 3.611388347828508
$ws.Range("I13").Value = 19.3215312479484
$ws.Range("K13").Value = 12.37378243093434
$ws.Range("M13").Value = 14.98615962537103
$ws.Range("O13").Value = 20.4761141815717
$ws.Range("B14").Value = 6.544554172909375
$ws.Range("D14").Value = 4.776237058912052
$ws.Range("E14").Value = 13.0827493921396
$ws.Range("F14").Value = 23.29112553356561
$ws.Range("G14").Value = 3.611740465434182
$ws.Range("I14").Value = 19.34316608525194
$ws.Range("K14").Value = 12.29467776611246
$ws.Range("M14").Value = 14.93085448273073
$ws.Range("O14").Value = 20.48037936744795
$ws.Range("B15").Value = 6.528254054625932
$ws.Range("D15").Value = 4.770505973446578
$ws.Range("E15").Value = 13.06250559934342
$ws.Range("F15").Value = 23.28430777256944
$ws.Range("G15").Value = 3.611957387910917
$ws.Range("I15").Value = 19.3564900247441
$ws.Range("K15").Value = 12.24593625034494
$ws.Range("M15").Value = 14.89692880693406
$ws.Range("O15").Value = 20.48310762100509
$ws.Range("B16").Value = 6.434206841105886
$ws.Range("D16").Value = 4.737424996523998
$ws.Range("E16").Value = 12.94733542256752
$ws.Range("F16").Value = 23.24749266642433
$ws.Range("G16").Value = 3.613219529208374
$ws.Range("I16").Value = 19.43394911879106
$ws.Range("K16").Value = 11.96223855217418
$ws.Range("M16").Value = 14.70180931104938
$ws.Range("O16").Value = 20.50051301298679
$ws.Range("B17").Value = 6.375970807400755
$ws.Range("D17").Value = 4.716922728897173
$ws.Range("E17").Value = 12.87746828406378
$ws.Range("F17").Value = 23.22692562675775
$ws.Range("G17").Value = 3.614010846823194
$ws.Range("I17").Value = 19.48245569489938
$ws.Range("K17").Value = 11.78432435986221
$ws.Range("M17").Value = 14.58154298540339
$ws.Range("O17").Value = 20.5127679168963
$ws.Range("B18").Value = 6.34228753943149
$ws.Range("D18").Value = 4.705055716361641
$ws.Range("E18").Value = 12.83758623647318
$ws.Range("F18").Value = 23.21583968756923
$ws.Range("G18").Value = 3.614472267706935
$ws.Range("I18").Value = 19.51071922522341
$ws.Range("K18").Value = 11.68058294564511
$ws.Range("M18").Value = 14.51218589378276
$ws.Range("O18").Value = 20.52039503183798
$ws.Range("B19").Value = 6.330851881505427
$ws.Range("D19").Value = 4.70102510570002
$ws.Range("E19").Value = 12.82413641569113
$ws.Range("F19").Value = 23.21221408885904
$ws.Range("G19").Value = 3.614629576277871
$ws.Range("I19").Value = 19.52035133250571
$ws.Range("K19").Value = 11.64521716543665
$ws.Range("M19").Value = 14.48867376742512
$ws.Range("O19").Value = 20.52307666467906
$ws.Range("B20").Value = 6.382189794469391
$ws.Range("D20").Value = 4.719112998314484
$ws.Range("E20").Value = 12.88487466301874
$ws.Range("F20").Value = 23.22903809992372
$ws.Range("G20").Value = 3.613925960531771
$ws.Range("I20").Value = 19.47725444933484
$ws.Range("K20").Value = 11.80340987984516
$ws.Range("M20").Value = 14.59436506607196
$ws.Range("O20").Value = 20.51140346883421
$ws.Range("B21").Value = 6.55235885491016
$ws.Range("D21").Value = 4.778980988092102
$ws.Range("E21").Value = 13.09247114920734
$ws.Range("F21").Value = 23.29443463023107
$ws.Range("G21").Value = 3.611636781572104
$ws.Range("I21").Value = 19.33679641339974
$ws.Range("K21").Value = 12.31797273049052
$ws.Range("M21").Value = 14.94710942992551
$ws.Range("O21").Value = 20.47910246041952
$ws.Range("B22").Value = 6.661668530485692
$ws.Range("D22").Value = 4.81740459785054
$ws.Range("E22").Value = 13.2305357975114
$ws.Range("F22").Value = 23.34374802365007
$ws.Range("G22").Value = 3.610196224760203
$ws.Range("I22").Value = 19.2482225504912
$ws.Range("K22").Value = 12.64143222522826
$ws.Range("M22").Value = 15.17549184802334
$ws.Range("O22").Value = 20.46316733960413
$ws.Range("B23").Value = 6.603513218872272
$ws.Range("D23").Value = 4.79696341646914
$ws.Range("E23").Value = 13.15664492773568
$ws.Range("F23").Value = 23.31682983710714
$ws.Range("G23").Value = 3.610960009867447
$ws.Range("I23").Value = 19.29520191346994
$ws.Range("K23").Value = 12.46998266433651
$ws.Range("M23").Value = 15.05382181124023
$ws.Range("O23").Value = 20.47119756184275
$ws.Range("B24").Value = 6.379378818777501
$ws.Range("D24").Value = 4.718123026202552
$ws.Range("E24").Value = 12.88152534629026
$ws.Range("F24").Value = 23.22808075043083
$ws.Range("G24").Value = 3.61396431742985
$ws.Range("I24").Value = 19.4796047596222
$ws.Range("K24").Value = 11.79478584866639
$ws.Range("M24").Value = 14.58856886517137
$ws.Range("O24").Value = 20.51201852437556
$ws.Range("B25").Value = 6.130611213458269
$ws.Range("D25").Value = 4.630251282394996
$ws.Range("E25").Value = 12.59655355192881
$ws.Range("F25").Value = 23.16079429322947
$ws.Range("G25").Value = 3.617445099086462
$ws.Range("I25").Value = 19.69243345347992
$ws.Range("K25").Value = 11.01281820610098
$ws.Range("M25").Value = 14.08006721630864
$ws.Range("O25").Value = 20.57814499268383
